# Add five new rows (24-28) to the bottom of the table on the active sheet,
# each containing a single script filename in column A - mirroring the
# existing rows (e.g. row 23) that only populate column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "SCRIPT/P02P01A/us0202.ssb",
    "SCRIPT/P02P01A/us0402.ssb",
    "SCRIPT/P02P01A/us2002.ssb",
    "SCRIPT/P02P01A/us2005.ssb",
    "SCRIPT/P02P01A/us2008.ssb"
)

$startRow = 24
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("A$row").Value = $values[$i]
    $ws.Rows.Item($row).RowHeight = 43.2
}
